$d = $word.ActiveDocument

$d.Content.Find.Execute("47×87=4089", $true, $false, $false, $false, $false, $true, 1, $false, "64×42=2688", 2)
$d.Content.Find.Execute("43×86=3698", $true, $false, $false, $false, $false, $true, 1, $false, "71×82=5822", 2)
$d.Content.Find.Execute("77×83=6391", $true, $false, $false, $false, $false, $true, 1, $false, "81×80=6480", 2)
$d.Content.Find.Execute("23×29=667", $true, $false, $false, $false, $false, $true, 1, $false, "87×49=4263", 2)
$d.Content.Find.Execute("55×81=4455", $true, $false, $false, $false, $false, $true, 1, $false, "39×77=3003", 2)
$d.Content.Find.Execute("79×26=2054", $true, $false, $false, $false, $false, $true, 1, $false, "79×79=6241", 2)
$d.Content.Find.Execute("36×55=1980", $true, $false, $false, $false, $false, $true, 1, $false, "86×89=7654", 2)
$d.Content.Find.Execute("97×14=1358", $true, $false, $false, $false, $false, $true, 1, $false, "38×96=3648", 2)
$d.Content.Find.Execute("73×19=1387", $true, $false, $false, $false, $false, $true, 1, $false, "68×66=4488", 2)
$d.Content.Find.Execute("71×95=6745", $true, $false, $false, $false, $false, $true, 1, $false, "69×89=6141", 2)
$d.Content.Find.Execute("84×11=924", $true, $false, $false, $false, $false, $true, 1, $false, "76×75=5700", 2)
$d.Content.Find.Execute("24×68=1632", $true, $false, $false, $false, $false, $true, 1, $false, "43×46=1978", 2)
$d.Content.Find.Execute("43×43=1849", $true, $false, $false, $false, $false, $true, 1, $false, "41×55=2255", 2)
$d.Content.Find.Execute("65×23=1495", $true, $false, $false, $false, $false, $true, 1, $false, "47×16=752", 2)
$d.Content.Find.Execute("15×56=840", $true, $false, $false, $false, $false, $true, 1, $false, "57×73=4161", 2)
$d.Content.Find.Execute("45×23=1035", $true, $false, $false, $false, $false, $true, 1, $false, "63×57=3591", 2)
$d.Content.Find.Execute("61×55=3355", $true, $false, $false, $false, $false, $true, 1, $false, "19×77=1463", 2)
$d.Content.Find.Execute("29×47=1363", $true, $false, $false, $false, $false, $true, 1, $false, "82×65=5330", 2)
$d.Content.Find.Execute("34×41=1394", $true, $false, $false, $false, $false, $true, 1, $false, "96×28=2688", 2)
$d.Content.Find.Execute("88×78=6864", $true, $false, $false, $false, $false, $true, 1, $false, "45×44=1980", 2)
$d.Content.Find.Execute("53×95=5035", $true, $false, $false, $false, $false, $true, 1, $false, "44×45=1980", 2)
$d.Content.Find.Execute("28×78=2184", $true, $false, $false, $false, $false, $true, 1, $false, "68×28=1904", 2)
$d.Content.Find.Execute("95×42=3990", $true, $false, $false, $false, $false, $true, 1, $false, "56×64=3584", 2)
$d.Content.Find.Execute("84×13=1092", $true, $false, $false, $false, $false, $true, 1, $false, "59×86=5074", 2)
$d.Content.Find.Execute("63×44=2772", $true, $false, $false, $false, $false, $true, 1, $false, "28×12=336", 2)
